$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: 2021年 - full data row ---
# Copy formatting (style index) from the "year label" cell A6, then overwrite the value.
$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("A7").Value = "2021年"

$ws.Range("B7").Value = 100.6
$ws.Range("C7").Value = 101.9
$ws.Range("D7").Value = 102.2
$ws.Range("E7").Value = 101.2
$ws.Range("F7").Value = 102.1
$ws.Range("G7").Value = 102.6
$ws.Range("H7").Value = 101.5
$ws.Range("I7").Value = 101.5
$ws.Range("J7").Value = 101.4

# --- Row 8: 2022年 - only C8 populated, rest present but blank ---
$ws.Range("A6").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "2022年"

$ws.Range("C8").Value = 101.8

# Materialize the remaining cells of row 8 as blank (present, but empty) cells,
# mirroring an untouched blank cell elsewhere on the sheet.
$ws.Range("Z100").Copy($ws.Range("B8"))
$ws.Range("Z100").Copy($ws.Range("D8:J8"))
